$wb = $excel.ActiveWorkbook

$sheetData = @{}

$sheetData["ALC"] = @{
  "H38" = 610.6667
  "I38" = 358.14285
  "J38" = 1494.5
  "K38" = 1074.42855
  "L38" = 4483.5
  "M38" = -702.4285500000001
  "N38" = -5227.5
  "H39" = 617.44446
  "I39" = 685.25
  "K39" = 2055.75
  "M39" = -1759.75
  "H41" = 1672.2858
  "I41" = 1716.7693
  "J41" = 1600
  "K41" = 1716.7693
  "L41" = 1600
  "M41" = -1276.7693
  "N41" = -2480
  "H70" = 1954.7142
  "I70" = 2240
  "K70" = 6720
  "M70" = -6450
  "H73" = 1954.7142
  "I73" = 2240
  "K73" = 6720
  "M73" = -5784
  "H98" = 6327.6924
  "I98" = 6605
  "K98" = 6605
  "M98" = -5107
  "H112" = 1919.7906
  "I112" = 784.2857
  "J112" = 2140.5833
  "K112" = 2352.8571
  "L112" = 6421.749899999999
  "M112" = -1244.8571
  "N112" = -8637.749899999999
  "H113" = 2771.8462
  "I113" = 2694.875
  "J113" = 2895
  "K113" = 2694.875
  "L113" = 2895
  "M113" = 559.125
  "N113" = -9403
  "H116" = 1992.4736
  "I116" = 1654.25
  "J116" = 2572.2856
  "K116" = 1654.25
  "L116" = 2572.2856
  "M116" = 1787.75
  "N116" = -9456.285599999999
  "H122" = 6327.6924
  "I122" = 6605
  "K122" = 19815
  "M122" = -17365
  "H132" = 4906503
  "I132" = 6292011.5
  "K132" = 18876034.5
  "M132" = -18873504.5
  "H141" = 549.1
  "I141" = 547.79486
  "J141" = 600
  "K141" = 1643.38458
  "L141" = 1800
  "M141" = 3536.61542
  "N141" = -12160
}

$sheetData["ARM"] = @{
  "H32" = 4583.9683
  "I32" = 4196.691
  "K32" = 4196.691
  "M32" = -3909.691
  "H97" = 394.53845
  "I97" = 410.75
  "J97" = 200
  "K97" = 410.75
  "L97" = 200
  "M97" = 85.25
  "N97" = -1192
  "H122" = 3234.875
  "I122" = 2903.4285
  "J122" = 5555
  "K122" = 8710.2855
  "L122" = 16665
  "M122" = -6260.2855
  "N122" = -21565
  "H131" = 49470
  "J131" = 49470
  "L131" = 49470
  "N131" = -59550
}

$sheetData["BSM"] = @{
  "H86" = 3693
  "I86" = 4391.2
  "J86" = 1947.5
  "K86" = 4391.2
  "L86" = 1947.5
  "M86" = -3268.2
  "N86" = -4193.5
  "H89" = 3693
  "I89" = 4391.2
  "J89" = 1947.5
  "K89" = 21956
  "L89" = 9737.5
  "M89" = -16340
  "N89" = -20969.5
  "H105" = 48091556
  "I105" = 56106560
  "J105" = 1537
  "K105" = 56106560
  "L105" = 1537
  "M105" = -56104813
  "N105" = -5031
  "H107" = 1810.9412
  "I107" = 1498.1428
  "K107" = 1498.1428
  "M107" = 421.8571999999999
}

$sheetData["CRP"] = @{
  "H16" = 66667868
  "I16" = 76924230
  "K16" = 76924230
  "M16" = -76923943
  "H58" = 833.2941
  "I58" = 738.0741
  "J58" = 1200.5714
  "K58" = 738.0741
  "L58" = 1200.5714
  "M58" = -535.0741
  "N58" = -1606.5714
  "H99" = 2260
  "I99" = 2212
  "K99" = 2212
  "M99" = -714
  "H113" = 66667868
  "I113" = 76924230
  "K113" = 76924230
  "M113" = -76922060
  "H122" = 1245.8
  "I122" = 1123.3077
  "J122" = 1473.2858
  "K122" = 3369.9231
  "L122" = 4419.857400000001
  "M122" = -919.9231
  "N122" = -9319.857400000001
  "H126" = 2260
  "I126" = 2212
  "K126" = 6636
  "M126" = -4166
  "H134" = 18519690
  "I134" = 1079.2632
  "J134" = 62501390
  "K134" = 3237.7896
  "L134" = 187504170
  "M134" = -702.7896000000001
  "N134" = -187509240
  "H136" = 833.2941
  "I136" = 738.0741
  "J136" = 1200.5714
  "K136" = 2214.2223
  "L136" = 3601.7142
  "M136" = 335.7776999999996
  "N136" = -8701.7142
}

$sheetData["CUL"] = @{
  "H131" = 21740486
  "J131" = 1443.1708
  "L131" = 4329.512400000001
  "N131" = -14409.5124
}

$sheetData["GSM"] = @{
  "H102" = 1702.5834
  "I102" = 1475.8572
  "K102" = 1475.8572
  "M102" = 146.1428000000001
  "H122" = 4779.8
  "I122" = 4633.3335
  "K122" = 13900.0005
  "M122" = -11450.0005
  "H132" = 1934.1428
  "I132" = 1821.9474
  "K132" = 5465.8422
  "M132" = -2935.8422
}

$sheetData["LTW"] = @{
  "H100" = 1668.5834
  "I100" = 1546.1428
  "K100" = 1546.1428
  "M100" = -1005.1428
  "H136" = 1067.875
  "I136" = 982.64703
  "J136" = 1550.8334
  "K136" = 2947.94109
  "L136" = 4652.5002
  "M136" = -397.9410899999998
  "N136" = -9752.5002
}

$sheetData["WVR"] = @{
  "H130" = 37940
  "J130" = 37940
  "L130" = 37940
  "N130" = -47980
  "H132" = 7842.1
  "I132" = 10630.143
  "K132" = 31890.429
  "M132" = -29360.429
  "H136" = 635.5294
  "I136" = 451.37036
  "K136" = 1354.11108
  "M136" = 1195.88892
  "H138" = 34339
  "J138" = 34339
  "L138" = 34339
  "N138" = -44619
}

$totalUpdates = 0
foreach ($sheetName in $sheetData.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)
  $cellMap = $sheetData[$sheetName]
  foreach ($cellRef in $cellMap.Keys) {
    $ws.Range($cellRef).Value = $cellMap[$cellRef]
    $totalUpdates = $totalUpdates + 1
  }
}

Write-Host "Applied $totalUpdates cell updates across $($sheetData.Keys.Count) sheets"